# Updated cryptos list on Sun Jul  7 15:29:23 UTC 2024 with GitHub Actions
#
# Applies updated Price (column D) and Volume(1h) (column E) values to the
# cryptocurrency table on the active worksheet. Values are written as plain
# text (matching the source data, which uses locale-style "." thousand
# separators and padded percentage strings), so for any value that Excel's
# COM layer would otherwise auto-convert to a number, we temporarily force
# the cell to Text format, assign the literal string, then restore the
# cell's original style/format so no other formatting is disturbed.

function Set-TextValue {
    param($cell, [string]$text)

    # Values such as "1.00", "7.33", "0.999" look like plain numbers and
    # Excel would silently coerce them (and mangle formatting / precision)
    # if assigned directly. Values containing two "." separators (e.g.
    # "56.886.80") or left-over whitespace/percent signs are never
    # misinterpreted, so they can be assigned as-is.
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $originalStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = $originalStyle
    } else {
        $cell.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "56.886.80"
Set-TextValue $ws.Range("E2") "  -0.49%  "
Set-TextValue $ws.Range("D3") "2.970.18"
Set-TextValue $ws.Range("E3") "  -1.47%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "499.82"
Set-TextValue $ws.Range("E5") "  -3.29%  "
Set-TextValue $ws.Range("D6") "137.89"
Set-TextValue $ws.Range("E6") "  -1.42%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.09%  "
Set-TextValue $ws.Range("E8") "  -2.01%  "
Set-TextValue $ws.Range("D9") "7.33"
Set-TextValue $ws.Range("E9") "  -3.31%  "
Set-TextValue $ws.Range("D10") "0.107"
Set-TextValue $ws.Range("E10") "  -2.45%  "
Set-TextValue $ws.Range("E11") "  -0.45%  "
Set-TextValue $ws.Range("D12") "3.473.39"
Set-TextValue $ws.Range("E13") "  -1.73%  "
Set-TextValue $ws.Range("D14") "25.86"
Set-TextValue $ws.Range("E14") "  -0.24%  "
Set-TextValue $ws.Range("D15") "0.0000159"
Set-TextValue $ws.Range("E15") "  -0.61%  "
Set-TextValue $ws.Range("D16") "56.950.32"
Set-TextValue $ws.Range("E16") "  -0.34%  "
Set-TextValue $ws.Range("E17") "  +1.16%  "
Set-TextValue $ws.Range("D18") "2.970.94"
Set-TextValue $ws.Range("E18") "  -1.55%  "
Set-TextValue $ws.Range("E19") "  -0.55%  "
Set-TextValue $ws.Range("D20") "7.82"
Set-TextValue $ws.Range("E20") "  -1.20%  "
Set-TextValue $ws.Range("D21") "319.41"
Set-TextValue $ws.Range("E21") "  -3.19%  "
Set-TextValue $ws.Range("D22") "0.999"
Set-TextValue $ws.Range("E22") "  -0.11%  "
Set-TextValue $ws.Range("E23") "  -0.77%  "
Set-TextValue $ws.Range("D24") "0.485"
Set-TextValue $ws.Range("E24") "  -0.21%  "
Set-TextValue $ws.Range("D25") "63.14"
Set-TextValue $ws.Range("E25") "  -1.08%  "
Set-TextValue $ws.Range("E26") "  -0.02%  "
Set-TextValue $ws.Range("E27") "  -5.18%  "
Set-TextValue $ws.Range("E28") "  -3.14%  "
Set-TextValue $ws.Range("D29") "6.52"
Set-TextValue $ws.Range("E29") "  -2.83%  "
Set-TextValue $ws.Range("D30") "7.07"
Set-TextValue $ws.Range("E30") "  -1.50%  "
Set-TextValue $ws.Range("E31") "  -3.20%  "
Set-TextValue $ws.Range("E32") "  -5.48%  "
Set-TextValue $ws.Range("D33") "20.12"
Set-TextValue $ws.Range("E33") "  -2.94%  "
Set-TextValue $ws.Range("D34") "154.69"
Set-TextValue $ws.Range("E34") "  -2.03%  "
Set-TextValue $ws.Range("D35") "4.60"
Set-TextValue $ws.Range("E35") "  -0.43%  "
Set-TextValue $ws.Range("D36") "5.74"
Set-TextValue $ws.Range("E36") "  -0.23%  "
Set-TextValue $ws.Range("E37") "  -3.35%  "
Set-TextValue $ws.Range("D38") "24.15"
Set-TextValue $ws.Range("E38") "  -0.45%  "
Set-TextValue $ws.Range("E39") "  -2.23%  "
Set-TextValue $ws.Range("D40") "37.63"
Set-TextValue $ws.Range("D41") "2.999.84"
Set-TextValue $ws.Range("E41") "  -1.51%  "
Set-TextValue $ws.Range("E42") "  -0.08%  "
Set-TextValue $ws.Range("D43") "3.73"
Set-TextValue $ws.Range("E43") "  +0.01%  "
Set-TextValue $ws.Range("E44") "  -1.78%  "
Set-TextValue $ws.Range("D45") "2.201.97"
Set-TextValue $ws.Range("E45") "  -4.14%  "
Set-TextValue $ws.Range("E46") "  -3.32%  "
Set-TextValue $ws.Range("D47") "0.946"
Set-TextValue $ws.Range("E47") "  -6.24%  "
Set-TextValue $ws.Range("D48") "5.93"
Set-TextValue $ws.Range("E48") "  +0.78%  "
Set-TextValue $ws.Range("D49") "0.0235"
Set-TextValue $ws.Range("E49") "  -2.75%  "
Set-TextValue $ws.Range("D50") "19.21"
Set-TextValue $ws.Range("E50") "  -0.96%  "
Set-TextValue $ws.Range("E51") "  -10.48%  "
